$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 1.3
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 1.83
$ws.Range("K5").Value = 2.25
$ws.Range("L5").Value = 11
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67
$ws.Range("U5").Value = 2.75
$ws.Range("V5").Value = 1.4
$ws.Range("X5").Value = 5
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 7.5
$ws.Range("AB5").Value = 41
$ws.Range("AD5").Value = 9.5
$ws.Range("AE5").Value = 34
$ws.Range("AH5").Value = 21
$ws.Range("AK5").Value = 201
$ws.Range("AL5").Value = 101
$ws.Range("AM5").Value = 126
$ws.Range("AN5").Value = 3
$ws.Range("AO5").Value = 6.5
$ws.Range("AQ5").Value = 19
$ws.Range("AS5").Value = 301
$ws.Range("AV5").Value = 101
$ws.Range("AW5").Value = 11
$ws.Range("AX5").Value = 51

# Row 14 updates
$ws.Range("N14").Value = 9

# Row 15 updates
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11
